$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - duplicate of Milan Raut's record (row 2) but with a new fill/border style
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Milan Raut"
$ws.Range("C10").Value = 23
$ws.Range("D10").Value = "M"
$ws.Range("E10").Value = "A"
$ws.Range("F10").Value = 3.8
$ws.Range("G10").Value = "Dang"
$ws.Range("A10:G10").Borders.LineStyle = 1
$ws.Range("A10:G10").Interior.Pattern = 1

# Row 11 - duplicate of Amir Shapkota's record (row 3)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Amir Shapkota"
$ws.Range("C11").Value = 22
$ws.Range("D11").Value = "M"
$ws.Range("E11").Value = "A"
$ws.Range("F11").Value = 3.9
$ws.Range("G11").Value = "Kathmandu"
$ws.Range("A11:G11").Borders.LineStyle = 1

# Row 12 - new record, Sailesh Shapkota
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Sailesh Shapkota"
$ws.Range("D12").Value = "M"
$ws.Range("E12").Value = "B"
$ws.Range("F12").Value = 2.7
$ws.Range("G12").Value = "Kathmandu"
$ws.Range("A12:G12").Borders.LineStyle = 1

# Row 13 - new record, Pratik Chaudhary
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Pratik Chaudhary"
$ws.Range("C13").Value = "Twenty"
$ws.Range("D13").Value = "M"
$ws.Range("E13").Value = "A"
$ws.Range("F13").Value = 3.4
$ws.Range("G13").Value = "Kapilbastu"
$ws.Range("A13:G13").Borders.LineStyle = 1

$ws.Range("G13").Select()
